$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the price/volume columns remain plain text (they are stored as
# text/inlineStr in the source data, e.g. "1.00", "183.33", "74.864.20")
# so Excel does not silently coerce the new values into numbers.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '74.864.20'
$ws.Range('E2').Value = '  +0.43%  '
$ws.Range('D3').Value = '2.812.12'
$ws.Range('E3').Value = '  +8.23%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '187.61'
$ws.Range('E5').Value = '  +0.30%  '
$ws.Range('D6').Value = '597.04'
$ws.Range('E6').Value = '  +1.76%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '0.556'
$ws.Range('E8').Value = '  +3.32%  '
$ws.Range('D9').Value = '0.192'
$ws.Range('E9').Value = '  -7.21%  '
$ws.Range('D10').Value = '2.810.90'
$ws.Range('E10').Value = '  +8.33%  '
$ws.Range('E11').Value = '  -1.04%  '
$ws.Range('E12').Value = '  +1.73%  '
$ws.Range('D13').Value = '4.84'
$ws.Range('E13').Value = '  +0.79%  '
$ws.Range('D14').Value = '3.327.10'
$ws.Range('E14').Value = '  +8.27%  '
$ws.Range('D15').Value = '74.873.65'
$ws.Range('E15').Value = '  +0.44%  '
$ws.Range('E16').Value = '  +2.20%  '
$ws.Range('D17').Value = '0.0000186'
$ws.Range('E17').Value = '  -3.44%  '
$ws.Range('D18').Value = '2.810.85'
$ws.Range('E18').Value = '  +8.23%  '
$ws.Range('D19').Value = '8.94'
$ws.Range('E19').Value = '  -2.57%  '
$ws.Range('D20').Value = '12.33'
$ws.Range('E20').Value = '  +4.04%  '
$ws.Range('D21').Value = '374.27'
$ws.Range('E21').Value = '  -0.55%  '
$ws.Range('D22').Value = '2.24'
$ws.Range('E22').Value = '  -1.80%  '
$ws.Range('E23').Value = '  -0.31%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').Value = '70.49'
$ws.Range('E25').Value = '  +0.22%  '
$ws.Range('D26').Value = '2.955.45'
$ws.Range('E26').Value = '  +8.26%  '
$ws.Range('E27').Value = '  -0.94%  '
$ws.Range('E28').Value = '  +1.92%  '
$ws.Range('E29').Value = '  +7.33%  '
$ws.Range('E30').Value = '  +0.10%  '
$ws.Range('D31').Value = '515.23'
$ws.Range('E31').Value = '  +0.94%  '
$ws.Range('E32').Value = '  -1.99%  '
$ws.Range('D33').Value = '7.85'
$ws.Range('E33').Value = '  -1.89%  '
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '162.85'
$ws.Range('E36').Value = '  +1.86%  '
$ws.Range('D37').Value = '20.07'
$ws.Range('E37').Value = '  +4.19%  '
$ws.Range('D38').Value = '0.119'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('D39').Value = '19.31'
$ws.Range('E39').Value = '  -0.36%  '
$ws.Range('B40').Value = 'Aave'
$ws.Range('C40').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D40').Value = '183.33'
$ws.Range('E40').Value = '  +17.00%  '
$ws.Range('B41').Value = 'USDe'
$ws.Range('C41').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D41').Value = '1.00'
$ws.Range('E41').Value = '  -0.02%  '
$ws.Range('E42').Value = '  +1.59%  '
$ws.Range('E43').Value = '  +3.07%  '
$ws.Range('D44').Value = '1.68'
$ws.Range('E44').Value = '  -1.21%  '
$ws.Range('E45').Value = '  +3.54%  '
$ws.Range('D46').Value = '39.75'
$ws.Range('E46').Value = '  +2.50%  '
$ws.Range('E47').Value = '  -3.23%  '
$ws.Range('E48').Value = '  -1.88%  '
$ws.Range('D49').Value = '0.566'
$ws.Range('E49').Value = '  +7.77%  '
$ws.Range('D50').Value = '3.73'
$ws.Range('E50').Value = '  +2.38%  '
$ws.Range('D51').Value = '0.612'
$ws.Range('E51').Value = '  +4.53%  '
